$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 11 (Q1/Q2/Q3 for student #8) - set to zero
$ws.Range("E11").Value = 0
$ws.Range("F11").Value = 0
$ws.Range("G11").Value = 0

# Row 12 (Q1/Q2/Q3 for student #9)
$ws.Range("E12").Value = 30
$ws.Range("F12").Value = 60
$ws.Range("G12").Value = 0
$ws.Range("M12").Value = "5.1: not exist overwrite and vircual and program trminate every input, not good output. 5.2: add func don’t work"

# Row 13 (Q1/Q2/Q3 for student #10)
$ws.Range("E13").Value = 100
$ws.Range("F13").Value = 85
$ws.Range("G13").Value = 10
$ws.Range("M13").Value = "5.2: rule is incorrect but have this message"

# Update selection to reflect the last active cell in the session
$ws.Range("H15").Select()
